$wb = $excel.ActiveWorkbook

# Sheet "2025"
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 39063.99109145206
$ws.Range("C2").Value = 483537.6274462014
$ws.Range("F2").Value = 94331.34471502228
$ws.Range("H2").Value = 25342.77928792104
$ws.Range("N2").Value = 23638.06126801545
$ws.Range("O2").Value = 19940.13531829346

# Sheet "2030"
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 30846.52922536713
$ws.Range("C2").Value = 1495599.874611417
$ws.Range("F2").Value = 70193.79982138964
$ws.Range("H2").Value = 56602.42752520426
$ws.Range("N2").Value = 51649.16401227913
$ws.Range("O2").Value = 42574.77934331147

# Sheet "2035"
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 242452.4252219552
$ws.Range("C2").Value = 943335.270081223
$ws.Range("F2").Value = 1425.925979620855
$ws.Range("H2").Value = 39373.98526588717
$ws.Range("N2").Value = 53308.16490721726
$ws.Range("O2").Value = 30023.09380555204

# Sheet "2040"
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 11578.49752443177

# Sheet "2045"
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 76705.58894163162
$ws.Range("B2").Value = 1930.947398408091
$ws.Range("N2").Value = 28147.3462746636
$ws.Range("O2").Value = 8312.661449003012
